{"js": "// Update the date line and the 25 division-fact cells in the practice\n// table. Each old value is unique in the document, and the pairs are\n// applied in the same order as the source diff, so a simple sequential\n// find & replace (matching the exact, case-sensitive text of each run)\n// reproduces the target edit without any cross-talk between pairs.\nconst replacements = [\n  [\"2023-09-05 Tuesday\", \"2023-09-06 Wednesday\"],\n  [\"78\u00f78=\", \"86\u00f73=\"],\n  [\"52\u00f75=\", \"27\u00f75=\"],\n  [\"68\u00f75=\", \"83\u00f79=\"],\n  [\"85\u00f74=\", \"67\u00f75=\"],\n  [\"31\u00f78=\", \"55\u00f72=\"],\n  [\"40\u00f79=\", \"87\u00f76=\"],\n  [\"21\u00f76=\", \"90\u00f78=\"],\n  [\"44\u00f78=\", \"77\u00f79=\"],\n  [\"41\u00f72=\", \"86\u00f79=\"],\n  [\"54\u00f73=\", \"51\u00f77=\"],\n  [\"82\u00f74=\", \"28\u00f77=\"],\n  [\"82\u00f79=\", \"85\u00f74=\"],\n  [\"82\u00f73=\", \"37\u00f79=\"],\n  [\"93\u00f72=\", \"67\u00f79=\"],\n  [\"78\u00f77=\", \"22\u00f75=\"],\n  [\"19\u00f78=\", \"92\u00f76=\"],\n  [\"38\u00f79=\", \"42\u00f75=\"],\n  [\"70\u00f74=\", \"35\u00f75=\"],\n  [\"35\u00f77=\", \"30\u00f77=\"],\n  [\"17\u00f79=\", \"42\u00f75=\"],\n  [\"53\u00f74=\", \"59\u00f74=\"],\n  [\"54\u00f74=\", \"57\u00f76=\"],\n  [\"76\u00f73=\", \"56\u00f77=\"],\n  [\"96\u00f74=\", \"28\u00f75=\"],\n  [\"99\u00f74=\", \"34\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division-fact cells in the practice\n# table. Each old value is unique in the document, and the pairs are\n# applied in the same order as the source diff, so a simple sequential\n# Find/Replace (wdReplaceAll, case-sensitive, whole text of each run)\n# reproduces the target edit without any cross-talk between pairs.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-05 Tuesday\", \"2023-09-06 Wednesday\"),\n    @(\"78\u00f78=\", \"86\u00f73=\"),\n    @(\"52\u00f75=\", \"27\u00f75=\"),\n    @(\"68\u00f75=\", \"83\u00f79=\"),\n    @(\"85\u00f74=\", \"67\u00f75=\"),\n    @(\"31\u00f78=\", \"55\u00f72=\"),\n    @(\"40\u00f79=\", \"87\u00f76=\"),\n    @(\"21\u00f76=\", \"90\u00f78=\"),\n    @(\"44\u00f78=\", \"77\u00f79=\"),\n    @(\"41\u00f72=\", \"86\u00f79=\"),\n    @(\"54\u00f73=\", \"51\u00f77=\"),\n    @(\"82\u00f74=\", \"28\u00f77=\"),\n    @(\"82\u00f79=\", \"85\u00f74=\"),\n    @(\"82\u00f73=\", \"37\u00f79=\"),\n    @(\"93\u00f72=\", \"67\u00f79=\"),\n    @(\"78\u00f77=\", \"22\u00f75=\"),\n    @(\"19\u00f78=\", \"92\u00f76=\"),\n    @(\"38\u00f79=\", \"42\u00f75=\"),\n    @(\"70\u00f74=\", \"35\u00f75=\"),\n    @(\"35\u00f77=\", \"30\u00f77=\"),\n    @(\"17\u00f79=\", \"42\u00f75=\"),\n    @(\"53\u00f74=\", \"59\u00f74=\"),\n    @(\"54\u00f74=\", \"57\u00f76=\"),\n    @(\"76\u00f73=\", \"56\u00f77=\"),\n    @(\"96\u00f74=\", \"28\u00f75=\"),\n    @(\"99\u00f74=\", \"34\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
